# Fix scrambled roster rows: player bio data (No., Player, Pos, Ht, Wt,
# Birth Date, Unnamed: 6, Exp, College, bbref url) for these four adjacent
# row pairs was shifted by one row; swap each pair back into alignment.
# The row-index column (A) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowPair($Row1, $Row2, $Columns) {
    # Snapshot both rows first so writes to Row1 don't clobber Row2's reads.
    $vals1 = @{}
    $vals2 = @{}
    foreach ($col in $Columns) {
        $vals1[$col] = $ws.Range($col + $Row1).Value2
        $vals2[$col] = $ws.Range($col + $Row2).Value2
    }

    foreach ($col in $Columns) {
        $v1 = $vals1[$col]
        $v2 = $vals2[$col]

        if ($null -eq $v2) {
            $ws.Range($col + $Row1).ClearContents()
        } else {
            $ws.Range($col + $Row1).Value = $v2
        }

        if ($null -eq $v1) {
            $ws.Range($col + $Row2).ClearContents()
        } else {
            $ws.Range($col + $Row2).Value = $v1
        }
    }
}

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

Swap-RowPair 4  5  $cols
Swap-RowPair 6  7  $cols
Swap-RowPair 8  9  $cols
Swap-RowPair 10 11 $cols
